$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows for B/C pairs within each year block (values, columns A:E)
$swapPairs = @(3,7,11,15)
foreach ($r1 in $swapPairs) {
    $r2 = $r1 + 1
    $range1 = $ws.Range("A$r1`:E$r1")
    $range2 = $ws.Range("A$r2`:E$r2")
    $vals1 = $range1.Value()
    $vals2 = $range2.Value()
    $range1.Value = $vals2
    $range2.Value = $vals1
}

# Delete columns F and G entirely (shifting nothing, just removing content & dimension shrinks)
$ws.Range("F1:G17").Delete()
